$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.208766341209412
$ws.Range("B1").Value = 1.550592541694641
$ws.Range("C1").Value = 6.99921989440918
$ws.Range("D1").Value = 2.193573236465454
$ws.Range("E1").Value = 1.172656059265137
